# current_sense_bom.xlsx update
# - Two DigiKey part-number corrections + one device rename (MAX9923H -> MAX9922)
# - A formatting pass applied to the data rows (2:23), which (re)creates a
#   dedicated cell style for the generic/no-border cells while leaving the
#   special left-aligned numeric cells (B9:B15) untouched.
# - Selection moved to C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part-number / value corrections -------------------------------------

# Row 2 (0402 capacitor, 0.1uF, C1): DigiKey part number corrected.
$ws.Range("G2").Value = "445-1304-1-ND"

# Row 22 (U1): MAX9923H -> MAX9922, plus its DigiKey part number.
$ws.Range("B22").Value = "MAX9922"
$ws.Range("C22").Value = "MAX9922"
$ws.Range("G22").Value = "MAX9922EUB+TCT-ND "

# --- Formatting pass over the data rows (2:23) ----------------------------
# Re-apply the "Normal" cell style to every populated cell that isn't part
# of the special left-aligned numeric block (B9:B15), mirroring the
# formatting refresh picked up in the saved workbook.

$ws.Range("A2:G2").Style = "Normal"
$ws.Range("A3:G3").Style = "Normal"
$ws.Range("A4:G4").Style = "Normal"
$ws.Range("A5:I5").Style = "Normal"

$ws.Range("A6").Style = "Normal"
$ws.Range("C6:F6").Style = "Normal"
$ws.Range("A7").Style = "Normal"
$ws.Range("C7:F7").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("C8:F8").Style = "Normal"

$ws.Range("A9").Style = "Normal"
$ws.Range("C9:G9").Style = "Normal"
$ws.Range("A10").Style = "Normal"
$ws.Range("C10:G10").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("C11:G11").Style = "Normal"
$ws.Range("A12").Style = "Normal"
$ws.Range("C12:G12").Style = "Normal"
$ws.Range("A13").Style = "Normal"
$ws.Range("C13:G13").Style = "Normal"
$ws.Range("A14").Style = "Normal"
$ws.Range("C14:G14").Style = "Normal"
$ws.Range("A15").Style = "Normal"
$ws.Range("C15:G15").Style = "Normal"

$ws.Range("A16:G16").Style = "Normal"
$ws.Range("A17:G17").Style = "Normal"
$ws.Range("A18:G18").Style = "Normal"
$ws.Range("A19:G19").Style = "Normal"
$ws.Range("A20:G20").Style = "Normal"

$ws.Range("A21:F21").Style = "Normal"
$ws.Range("I21").Style = "Normal"

$ws.Range("A22:G22").Style = "Normal"
$ws.Range("A23:G23").Style = "Normal"

# --- Selection -------------------------------------------------------------
[void]$ws.Range("C17").Select()
